{"js": "// Rename the R objects `data` -> `athletes` and `lin.reg` -> `linReg`\n// throughout the \"Chp 13 Example 7\" R script walkthrough, matching the\n// commit that edited the RScript/annotated output for this example.\n//\n// We operate paragraph-by-paragraph and pick out the specific\n// occurrence(s) of the ambiguous token `data` that refer to the data-\n// frame variable (as opposed to the `data =` named argument, the\n// \"Reading in data\" heading, or the `data` that is just part of the\n// CSV's URL) so that only the intended identifier is renamed, while\n// run-level formatting (syntax-highlighting styles such as NormalTok /\n// FunctionTok / VerbatimChar) is preserved because we replace text only\n// inside the matched sub-range.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- Paragraph 6: \"data <- read.csv(file='...college_female_athletes.csv')\"\n// Rename only the leading variable name (first \"data\"); the \"data\" inside\n// the URL string must stay untouched.\n{\n  const p = paragraphs.items[6];\n  const matches = p.search(\"data\", { matchCase: true, matchWholeWord: true });\n  matches.load(\"items/text\");\n  await context.sync();\n  if (matches.items.length > 0) {\n    matches.items[0].insertText(\"athletes\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// --- Paragraph 8: \"lin.reg <- lm(TBW ~ HGT + BF + AGE, data = data)\\nlin.reg\"\n// Rename both \"lin.reg\" occurrences to \"linReg\", and rename only the\n// second \"data\" (the data-frame argument value) to \"athletes\" -- the\n// first \"data\" (the \"data =\" argument name) must stay untouched.\n{\n  const p = paragraphs.items[8];\n  const linregMatches = p.search(\"lin.reg\", { matchCase: true });\n  linregMatches.load(\"items/text\");\n  const dataMatches = p.search(\"data\", { matchCase: true, matchWholeWord: true });\n  dataMatches.load(\"items/text\");\n  await context.sync();\n\n  for (let i = 0; i < linregMatches.items.length; i++) {\n    linregMatches.items[i].insertText(\"linReg\", Word.InsertLocation.replace);\n  }\n  if (dataMatches.items.length > 1) {\n    dataMatches.items[1].insertText(\"athletes\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// --- Paragraph 9: \"## lm(formula = TBW ~ HGT + BF + AGE, data = data)\" (plus\n// the rest of the printed model-summary block). Rename only the trailing\n// \"data\" (the argument value) to \"athletes\"; \"data =\" stays untouched.\n{\n  const p = paragraphs.items[9];\n  const dataMatches = p.search(\"data\", { matchCase: true, matchWholeWord: true });\n  dataMatches.load(\"items/text\");\n  await context.sync();\n  if (dataMatches.items.length > 1) {\n    dataMatches.items[1].insertText(\"athletes\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// --- Paragraph 11: \"round(confint(lin.reg), 1)\"\n{\n  const p = paragraphs.items[11];\n  const linregMatches = p.search(\"lin.reg\", { matchCase: true });\n  linregMatches.load(\"items/text\");\n  await context.sync();\n  for (let i = 0; i < linregMatches.items.length; i++) {\n    linregMatches.items[i].insertText(\"linReg\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Rename the R objects `data` -> `athletes` and `lin.reg` -> `linReg`\n# throughout the \"Chp 13 Example 7\" R script walkthrough, matching the\n# commit that edited the RScript/annotated output for this example.\n#\n# Replacements are scoped to individual paragraphs (by 1-based COM\n# index) and target a specific occurrence of the ambiguous token\n# `data` so that only the data-frame *variable* is renamed, while the\n# `data =` named argument, the \"Reading in data\" heading, and the\n# `data` that is part of the CSV's URL are left untouched. Using\n# Range.Find (rather than rewriting whole paragraphs) keeps each run's\n# syntax-highlighting style (NormalTok / FunctionTok / VerbatimChar /\n# ...) intact.\n\nfunction Replace-NthMatch {\n    param($Doc, $ParaIndex, $SearchText, $ReplaceText, $Occurrence, $WholeWord)\n    $p = $Doc.Paragraphs($ParaIndex)\n    $pStart = $p.Range.Start\n    $pEnd = $p.Range.End\n\n    $r = $Doc.Range($pStart, $pEnd)\n    $r.Find.ClearFormatting()\n    $r.Find.Text = $SearchText\n    $r.Find.MatchWholeWord = $WholeWord\n    $r.Find.MatchCase = $true\n    $r.Find.Forward = $true\n    $r.Find.Wrap = 0\n\n    $count = 0\n    while ($r.Find.Execute()) {\n        if ($r.Start -ge $pEnd) { break }\n        $count++\n        if ($count -eq $Occurrence) {\n            $r.Text = $ReplaceText\n            return $true\n        }\n        $r.Start = $r.End\n        $r.End = $pEnd\n    }\n    return $false\n}\n\nfunction Replace-AllMatches {\n    param($Doc, $ParaIndex, $SearchText, $ReplaceText, $WholeWord)\n    $p = $Doc.Paragraphs($ParaIndex)\n    $pStart = $p.Range.Start\n    $pEnd = $p.Range.End\n\n    $r = $Doc.Range($pStart, $pEnd)\n    $r.Find.ClearFormatting()\n    $r.Find.Text = $SearchText\n    $r.Find.MatchWholeWord = $WholeWord\n    $r.Find.MatchCase = $true\n    $r.Find.Forward = $true\n    $r.Find.Wrap = 0\n\n    $lenDiff = $ReplaceText.Length - $SearchText.Length\n    while ($r.Find.Execute()) {\n        if ($r.Start -ge $pEnd) { break }\n        $r.Text = $ReplaceText\n        $pEnd = $pEnd + $lenDiff\n        $r.Start = $r.End\n        $r.End = $pEnd\n    }\n}\n\n$d = $word.ActiveDocument\n\n# Paragraph 7: \"data <- read.csv(file='...college_female_athletes.csv')\"\n# Rename only the leading variable name (1st \"data\"); the \"data\" inside\n# the URL string must stay untouched.\nReplace-NthMatch $d 7 \"data\" \"athletes\" 1 $true | Out-Null\n\n# Paragraph 9: \"lin.reg <- lm(TBW ~ HGT + BF + AGE, data = data)\" / \"lin.reg\"\n# Rename both \"lin.reg\" occurrences, and only the 2nd \"data\" (the\n# data-frame argument value) -- the 1st \"data\" is the \"data =\" argument\n# name and must stay untouched.\nReplace-AllMatches $d 9 \"lin.reg\" \"linReg\" $false\nReplace-NthMatch $d 9 \"data\" \"athletes\" 2 $true | Out-Null\n\n# Paragraph 10: \"## lm(formula = TBW ~ HGT + BF + AGE, data = data)\" (output)\n# Rename only the trailing \"data\" (the argument value); \"data =\" stays.\nReplace-NthMatch $d 10 \"data\" \"athletes\" 2 $true | Out-Null\n\n# Paragraph 12: \"round(confint(lin.reg), 1)\"\nReplace-AllMatches $d 12 \"lin.reg\" \"linReg\" $false\n"}
